# Update the "Final_Non-Linear_Wands_Data" sheet:
#  - New (cleaned up) cost values for column L (rows 2-60)
#  - Apply a "0" integer number format to the updated L cells
#  - Update the sheet's scroll position / active selection to the L2:L60 range

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

$lValues = @{
    2 = 9
    3 = 13
    4 = 19
    5 = 28
    6 = 40
    7 = 59
    8 = 87
    9 = 127
    10 = 186
    11 = 273
    12 = 399
    13 = 585
    14 = 856
    15 = 1254
    16 = 1837
    17 = 2691
    18 = 3941
    19 = 5772
    20 = 8454
    21 = 12382
    22 = 18135
    23 = 26561
    24 = 38902
    25 = 56977
    26 = 83451
    27 = 122226
    28 = 179017
    29 = 262195
    30 = 384022
    31 = 562454
    32 = 823793
    33 = 1206561
    34 = 1767179
    35 = 2588282
    36 = 3790904
    37 = 5552314
    38 = 8132146
    39 = 11910675
    40 = 17444863
    41 = 25550461
    42 = 37422253
    43 = 54810166
    44 = 80277216
    45 = 117577302
    46 = 172208537
    47 = 252223684
    48 = 369417149
    49 = 541063502
    50 = 792463789
    51 = 1160674956
    52 = 1699972128
    53 = 2000000000
    54 = 3646734750
    55 = 5341157231
    56 = 7822877869
    57 = 11457707664
    58 = 16781428411
    59 = 24578768089
    60 = 35999071473
}

foreach ($row in $lValues.Keys) {
    $ws.Cells.Item([int]$row, 12).Value = $lValues[$row]
}

$lRange = $ws.Range("L2:L60")
$lRange.NumberFormat = "0"

# Scroll the view so column I is the left-most visible column (best effort;
# mirrors the saved sheetView's topLeftCell), then select L2:L60 so it
# becomes the sheet's active selection (activeCell L2, sqref L2:L60).
$excel.ActiveWindow.ScrollColumn = 9
$excel.ActiveWindow.ScrollRow = 1

$lRange.Select()
